# cardData v2, create list
# Swap the "이동범위" (move range) and "효과" (effect) columns on the
# 커맨더 (commander) sheet: what used to be column D (move range) becomes
# column E, and what used to be column E (effect) becomes column D -
# both for the header row (row 1) and the first data row (row 2).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Read current contents (header + first data row) for columns D and E.
$d1 = $ws1.Range("D1").Text
$e1 = $ws1.Range("E1").Text
$d2 = $ws1.Range("D2").Text
$e2 = $ws1.Range("E2").Text

# Write them back swapped.
$ws1.Range("D1").Value = $e1
$ws1.Range("E1").Value = $d1
$ws1.Range("D2").Value = $e2
$ws1.Range("E2").Value = $d2

# Restore the selection / active sheet state left by the author.
$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate()
$ws3.Range("C4").Select()

$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Range("E14").Select()
